$d = $word.ActiveDocument

# Update the date line at the top of the document
[void]$d.Content.Find.Execute("2024-07-23 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-07-24 Wednesday", 2)

# Update the division problems in the table, cell by cell, so that
# replacement values which collide with other pre-existing source values
# (e.g. "13÷6=" is both an original value and a new value elsewhere)
# do not get double-replaced.
$tbl = $d.Tables.Item(1)

$edits = @(
    @{ Row = 1;  Col = 1; New = "17÷7=" },
    @{ Row = 1;  Col = 2; New = "72÷4=" },
    @{ Row = 1;  Col = 3; New = "99÷9=" },
    @{ Row = 1;  Col = 4; New = "18÷9=" },
    @{ Row = 1;  Col = 5; New = "87÷8=" },

    @{ Row = 5;  Col = 1; New = "85÷4=" },
    @{ Row = 5;  Col = 2; New = "76÷7=" },
    @{ Row = 5;  Col = 3; New = "97÷3=" },
    @{ Row = 5;  Col = 4; New = "94÷3=" },
    @{ Row = 5;  Col = 5; New = "28÷6=" },

    @{ Row = 9;  Col = 1; New = "76÷6=" },
    @{ Row = 9;  Col = 2; New = "99÷8=" },
    @{ Row = 9;  Col = 3; New = "58÷6=" },
    @{ Row = 9;  Col = 4; New = "99÷2=" },
    @{ Row = 9;  Col = 5; New = "41÷9=" },

    @{ Row = 13; Col = 1; New = "97÷3=" },
    @{ Row = 13; Col = 2; New = "10÷7=" },
    @{ Row = 13; Col = 3; New = "13÷6=" },
    @{ Row = 13; Col = 4; New = "71÷4=" },
    @{ Row = 13; Col = 5; New = "91÷7=" },

    @{ Row = 17; Col = 1; New = "65÷9=" },
    @{ Row = 17; Col = 2; New = "42÷3=" },
    @{ Row = 17; Col = 3; New = "63÷8=" },
    @{ Row = 17; Col = 4; New = "66÷3=" },
    @{ Row = 17; Col = 5; New = "74÷4=" }
)

foreach ($edit in $edits) {
    $cell = $tbl.Cell($edit.Row, $edit.Col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $edit.New
}
